$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colB = New-Object 'double[,]' 24,1
$colB[0,0] = 0.330636746856328
$colB[1,0] = 0.291116090231867
$colB[2,0] = 0.2668413665207368
$colB[3,0] = 0.2569474271382717
$colB[4,0] = 0.2553044507996844
$colB[5,0] = 0.2667079400673344
$colB[6,0] = 0.3170121197501032
$colB[7,0] = 0.4155734735245176
$colB[8,0] = 0.487922311229255
$colB[9,0] = 0.5208195530175317
$colB[10,0] = 0.533274434928785
$colB[11,0] = 0.530592176211087
$colB[12,0] = 0.5218442783456965
$colB[13,0] = 0.5164855865738502
$colB[14,0] = 0.4857720633422105
$colB[15,0] = 0.4669262220450321
$colB[16,0] = 0.4560852468404448
$colB[17,0] = 0.4524144626962538
$colB[18,0] = 0.4689325368490813
$colB[19,0] = 0.524413822173841
$colB[20,0] = 0.5606587290533867
$colB[21,0] = 0.5413157032594142
$colB[22,0] = 0.4680255013583121
$colB[23,0] = 0.3889203981700575
$ws.Range("B2:B25").Value = $colB

$colC = New-Object 'double[,]' 24,1
$colC[0,0] = 0.06440320679726597
$colC[1,0] = 0.05667629207825087
$colC[2,0] = 0.05190885488246977
$colC[3,0] = 0.049960301766518
$colC[4,0] = 0.04963639728343594
$colC[5,0] = 0.05188259940453577
$colC[6,0] = 0.06174376659832603
$colC[7,0] = 0.08089839570087065
$colC[8,0] = 0.09486110442105655
$colC[9,0] = 0.1011895636161455
$colC[10,0] = 0.1035826445485668
$colC[11,0] = 0.1030674020149718
$colC[12,0] = 0.1013865114993564
$colC[13,0] = 0.1003564779754811
$colC[14,0] = 0.09444705468180814
$colC[15,0] = 0.09081584193891956
$colC[16,0] = 0.08872507784583661
$colC[17,0] = 0.08801680631174236
$colC[18,0] = 0.0912026172284186
$colC[19,0] = 0.1018803214363118
$colC[20,0] = 0.1088391680290783
$colC[21,0] = 0.1051269072554248
$colC[22,0] = 0.09102776587369021
$colC[23,0] = 0.07573601649018258
$ws.Range("C2:C25").Value = $colC

$colD = New-Object 'double[,]' 24,1
$colD[0,0] = 0.02380209259330002
$colD[1,0] = 0.02184017161696516
$colD[2,0] = 0.02062559912396011
$colD[3,0] = 0.02012818275436956
$colD[4,0] = 0.0200454390303193
$colD[5,0] = 0.02061890074221395
$colD[6,0] = 0.0231277052973482
$colD[7,0] = 0.02796733855954159
$colD[8,0] = 0.03147285446512171
$colD[9,0] = 0.03305645241681532
$colD[10,0] = 0.03365449761152917
$colD[11,0] = 0.03352577089775366
$colD[12,0] = 0.03310568680674919
$colD[13,0] = 0.03284815996740065
$colD[14,0] = 0.03136913661330709
$colD[15,0] = 0.03045894059690113
$colD[16,0] = 0.02993437926851072
$colD[17,0] = 0.02975659417704435
$colD[18,0] = 0.03055594048008459
$colD[19,0] = 0.0332291201829662
$colD[20,0] = 0.03496668612122278
$colD[21,0] = 0.03404019661113722
$colD[22,0] = 0.03051209080929596
$colD[23,0] = 0.02666680388576026
$ws.Range("D2:D25").Value = $colD

$colE = New-Object 'double[,]' 24,1
$colE[0,0] = 0.4180825186005848
$colE[1,0] = 0.3648139447740135
$colE[2,0] = 0.3321876317311023
$colE[3,0] = 0.3189109501738869
$colE[4,0] = 0.3167074590422061
$colE[5,0] = 0.3320085036036176
$colE[6,0] = 0.3996978156532975
$colE[7,0] = 0.5331519221147261
$colE[8,0] = 0.6317487387890992
$colE[9,0] = 0.676746716428454
$colE[10,0] = 0.6938090645357846
$colE[11,0] = 0.6901333518449206
$colE[12,0] = 0.6781499839690923
$colE[13,0] = 0.6708128155604243
$colE[14,0] = 0.6288110875301101
$colE[15,0] = 0.6030828814736537
$colE[16,0] = 0.5882983224978489
$colE[17,0] = 0.5832948155519944
$colE[18,0] = 0.6058202711046619
$colE[19,0] = 0.6816691634749645
$colE[20,0] = 0.7313732378805753
$colE[21,0] = 0.7048325792938783
$colE[22,0] = 0.6045826758208079
$colE[23,0] = 0.4969609217057638
$ws.Range("E2:E25").Value = $colE

$colF = New-Object 'double[,]' 24,1
$colF[0,0] = 0.6122752868145653
$colF[1,0] = 0.6089813308479037
$colF[2,0] = 0.6073720781820597
$colF[3,0] = 0.6068201013650736
$colF[4,0] = 0.6067347125572908
$colF[5,0] = 0.6073642138654947
$colF[6,0] = 0.61105369264709
$colF[7,0] = 0.621573946837259
$colF[8,0] = 0.6313165380492407
$colF[9,0] = 0.6361882923576871
$colF[10,0] = 0.6380964993384737
$colF[11,0] = 0.63768271183217
$colF[12,0] = 0.6363440106950122
$colF[13,0] = 0.6355322755128157
$colF[14,0] = 0.6310070190307826
$colF[15,0] = 0.6283436613675448
$colF[16,0] = 0.6268531564255539
$colF[17,0] = 0.6263556013268925
$colF[18,0] = 0.6286228958501496
$colF[19,0] = 0.6367354984828921
$colF[20,0] = 0.6424070401718325
$colF[21,0] = 0.6393461783937227
$colF[22,0] = 0.6284965271862362
$colF[23,0] = 0.6183751400647424
$ws.Range("F2:F25").Value = $colF

$colH = New-Object 'double[,]' 24,1
$colH[0,0] = 0.07973214163530429
$colH[1,0] = 0.07973214163530429
$colH[2,0] = 0.07973214163530429
$colH[3,0] = 0.07973214163530429
$colH[4,0] = 0.07973214163530429
$colH[5,0] = 0.07973214163530429
$colH[6,0] = 0.07973214163530429
$colH[7,0] = 0.07973214163530429
$colH[8,0] = 0.07973214163530429
$colH[9,0] = 0.07973214163530429
$colH[10,0] = 0.07973214163530429
$colH[11,0] = 0.07973214163530429
$colH[12,0] = 0.07973214163530429
$colH[13,0] = 0.07973214163530429
$colH[14,0] = 0.07973214163530429
$colH[15,0] = 0.07973214163530429
$colH[16,0] = 0.07973214163530429
$colH[17,0] = 0.07973214163530429
$colH[18,0] = 0.07973214163530429
$colH[19,0] = 0.07973214163530429
$colH[20,0] = 0.07973214163530429
$colH[21,0] = 0.07973214163530429
$colH[22,0] = 0.07973214163530429
$colH[23,0] = 0.07973214163530429
$ws.Range("H2:H25").Value = $colH

$colI = New-Object 'double[,]' 24,1
$colI[0,0] = 0.4517866951497496
$colI[1,0] = 0.4552154285460013
$colI[2,0] = 0.4575912503754154
$colI[3,0] = 0.4586273773408429
$colI[4,0] = 0.4588035278624432
$colI[5,0] = 0.4576049488791263
$colI[6,0] = 0.4529127206142505
$colI[7,0] = 0.4458617095770414
$colI[8,0] = 0.4419975978579913
$colI[9,0] = 0.4405267024191915
$colI[10,0] = 0.4400110636896599
$colI[11,0] = 0.4401202747286774
$colI[12,0] = 0.440483451069138
$colI[13,0] = 0.4407112960116208
$colI[14,0] = 0.4420995062339941
$colI[15,0] = 0.4430246836351301
$colI[16,0] = 0.4435838203440596
$colI[17,0] = 0.4437777685057362
$colI[18,0] = 0.4429234018844532
$colI[19,0] = 0.4403756541427555
$colI[20,0] = 0.4389516634960913
$colI[21,0] = 0.439689579260147
$colI[22,0] = 0.4429691065181487
$colI[23,0] = 0.4475384298113809
$ws.Range("I2:I25").Value = $colI

$colK = New-Object 'double[,]' 24,1
$colK[0,0] = 0.3611763121005538
$colK[1,0] = 0.315793928048123
$colK[2,0] = 0.287886050283845
$colK[3,0] = 0.276502986616066
$colK[4,0] = 0.2746122230464607
$colK[5,0] = 0.2877325757725941
$colK[6,0] = 0.3455376061479285
$colK[7,0] = 0.4585399660348628
$colK[8,0] = 0.5413386479538644
$colK[9,0] = 0.5789558663396974
$colK[10,0] = 0.5931932769146044
$colK[11,0] = 0.5901273329321839
$colK[12,0] = 0.5801273379195493
$colK[13,0] = 0.574001072957202
$colK[14,0] = 0.5388792621645848
$colK[15,0] = 0.5173204619230489
$colK[16,0] = 0.5049158901433373
$colK[17,0] = 0.500715153156051
$colK[18,0] = 0.5196159042024817
$colK[19,0] = 0.5830647837674121
$colK[20,0] = 0.6244888668702515
$colK[21,0] = 0.6023841812906596
$colK[22,0] = 0.518578166269549
$colK[23,0] = 0.4280085604363819
$ws.Range("K2:K25").Value = $colK

$colN = New-Object 'double[,]' 24,1
$colN[0,0] = 1.189202072763623
$colN[1,0] = 1.203388477486316
$colN[2,0] = 1.212544879463685
$colN[3,0] = 1.216388262846273
$colN[4,0] = 1.21703322144032
$colN[5,0] = 1.212596258953303
$colN[6,0] = 1.19400097243277
$colN[7,0] = 1.161075244190114
$colN[8,0] = 1.139043985527588
$colN[9,0] = 1.129490597578853
$colN[10,0] = 1.125940445706279
$colN[11,0] = 1.126702030265738
$colN[12,0] = 1.12919717086247
$colN[13,0] = 1.130734311310594
$colN[14,0] = 1.139677763028086
$colN[15,0] = 1.14528441027444
$colN[16,0] = 1.148553319260373
$colN[17,0] = 1.149667689139507
$colN[18,0] = 1.144683006743454
$colN[19,0] = 1.128462454569897
$colN[20,0] = 1.118254929640381
$colN[21,0] = 1.123666833282691
$colN[22,0] = 1.144954759324472
$colN[23,0] = 1.169603452473941
$ws.Range("N2:N25").Value = $colN

$colO = New-Object 'double[,]' 24,1
$colO[0,0] = 2.053608104493634
$colO[1,0] = 2.05573438518465
$colO[2,0] = 2.058376677175673
$colO[3,0] = 2.059789110565859
$colO[4,0] = 2.060043904841592
$colO[5,0] = 2.058394367213054
$colO[6,0] = 2.054063522378442
$colO[7,0] = 2.05620070256208
$colO[8,0] = 2.064287731982063
$colO[9,0] = 2.069389979268919
$colO[10,0] = 2.071527350385111
$colO[11,0] = 2.071057891056029
$colO[12,0] = 2.069561705041878
$colO[13,0] = 2.068671996034027
$colO[14,0] = 2.06398298002199
$colO[15,0] = 2.0614714002987
$colO[16,0] = 2.06016075025002
$colO[17,0] = 2.05973997423294
$colO[18,0] = 2.06172489536047
$colO[19,0] = 2.069995595826555
$colO[20,0] = 2.076597662456351
$colO[21,0] = 2.072964327735122
$colO[22,0] = 2.061609875080592
$colO[23,0] = 2.054480762448691
$ws.Range("O2:O25").Value = $colO

Write-Host "Updated pl_mw values for 380 kV case"